$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "financing" columns/headers for agents (L = days_between_financing, M = financing_period)
$ws.Range("L1").Value = "days_between_financing"
$ws.Range("M1").Value = "financing_period"

# New per-agent values
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 10

$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 20

$ws.Range("L4").Value = 7
$ws.Range("M4").Value = 30

# Match the centered alignment style already used by the rest of the header/data block
$ws.Range("L1:M4").HorizontalAlignment = -4108
$ws.Range("L1:M4").VerticalAlignment = -4108

# Give the new columns an explicit (bestFit-like) width, same as the other labeled columns
$ws.Columns.Item(12).ColumnWidth = 144 / 7
$ws.Columns.Item(13).ColumnWidth = 95 / 7

# Move the active selection to the newly-populated M4 cell
$ws.Range("M4").Select()
